$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for specific rows as per repulled data / mean calculation fix
$ws.Range("F5").Value = -4
$ws.Range("F10").Value = -11
$ws.Range("F11").Value = -5
$ws.Range("F20").Value = -5
$ws.Range("F22").Value = -5
$ws.Range("F27").Value = 3
$ws.Range("F28").Value = -5
$ws.Range("F31").Value = -4
